# Add "Login with Facebook" feature — color key requirement bullets blue
# (0070C0), merge the split "Email" / ", Facebook " runs into one, and
# shift Word's internal "_GoBack" last-edit bookmark from the old
# "Knowledge (Dro|pdown)" split to a new split inside
# "ของบทความทั้งหมด..." (the net effect of the author's edits elsewhere
# in the doc).

$d = $word.ActiveDocument
$wdColorBlue0070C0 = 12611584   # RGB(0x00,0x70,0xC0) packed as BGR OLE color

# ------------------------------------------------------------------
# 1. "Member" paragraph -> whole paragraph (incl. paragraph mark) blue
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Member", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$memberPara = $rng.Paragraphs(1)
$memberPara.Range.Font.Color = $wdColorBlue0070C0

# ------------------------------------------------------------------
# 2. "...สามารถสมัครสมาชิกผ่าน Email, Facebook ...ได้" paragraph
#    a) merge "Email" + ", Facebook " runs into a single run via a
#       self-replace (Word coalesces identically-formatted adjacent
#       runs when the match is rewritten)
#    b) color the whole paragraph (incl. mark + every run) blue
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Email, Facebook ", $true, $false, $false, $false, $false, $true, 1, $false, "Email, Facebook ", 2)

$rng3 = $d.Content
$rng3.Find.Execute("สามารถสมัครสมาชิกผ่าน", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$memberFeaturePara = $rng3.Paragraphs(1)
$memberFeaturePara.Range.Font.Color = $wdColorBlue0070C0

# ------------------------------------------------------------------
# 3. "Knowledge (Dro" + "pdown " runs (currently split around the old
#    "_GoBack" bookmark) merge back into one run "Knowledge (Dropdown "
#    -- rewriting the text removes the stale bookmark automatically.
# ------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("Knowledge (Dropdown ", $true, $false, $false, $false, $false, $true, 1, $false, "Knowledge (Dropdown ", 2)

# ------------------------------------------------------------------
# 4. Relocate the "_GoBack" bookmark: split the run
#    "ของบทความทั้งหมด ตำแหน่งถัดจาก " after its first character "ข"
#    and drop a fresh zero-length "_GoBack" bookmark there (this also
#    removes the old bookmark, since Word keeps only one).
# ------------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.Execute("ของบทความทั้งหมด ตำแหน่งถัดจาก ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rng5.Start + 1
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "edit complete"
